# PAS-6576 Update "individual VIN retrieval" logic to use ENTRY DATE and VALID
# Apply the same data edits made to AddedVIN_CA_SELECT.xlsx:
#  - MODEL_TEXT (col F) on rows 3-5 becomes distinguishing VIN-upload markers
#  - ENTRYDATE (col AI) on row 2 bumped a year
#  - sheet view selection/scroll updated to where the author was working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data changes -----------------------------------------------------

# Row 2: ENTRYDATE 2000-01-01 -> 2001-01-01
$ws.Range("AI2").Value = 20010101

# Rows 3-5: MODEL_TEXT "Gt" -> distinct markers used by the updated
# individual-VIN-retrieval tests
$ws.Range("F3").Value = "invalidVIN"
$ws.Range("F4").Value = "SecondValid"
$ws.Range("F5").Value = "ThirdValid"

# --- view/selection changes --------------------------------------------

# Best-effort: real Excel stamps the absolute save path into
# xl/workbook.xml (x15ac:absPath) on save; that metadata reflects the
# machine the file was last saved from and isn't an author-controllable
# document property, so there's no COM surface to drive it deliberately.
try {
    $wb.Path = "C:\Users\gu1xkaz\IdeaProjects\pas-cvqaautomation\aaa-automation-tests\src\test\resources\uploadingfiles\vinUploadFiles\"
} catch {
}

# Move the selection to where the author ended up (also nudges the
# window's visible/scrolled range).
$ws.Range("AI15").Select()
